$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily entries (2025-10-22 -> serial 45952) added to the bottom of the log,
# rows 499-508. Columns: A Date, B Nom du joueur, C Volume, D Intensite,
# E Fatigue, F Douleur, G Localisation douleur, H Plaisir, I Charge (=C*D).

$data = @(
    @{ Row=499; Name="Maé Clavel";       C=70; D=5; E=5; F=0; G=$null;           H=4  },
    @{ Row=500; Name="Yoann Martelat";   C=70; D=4; E=5; F=4; G="Genou";         H=7  },
    @{ Row=501; Name="Ilyes Boughanmi";  C=70; D=6; E=5; F=5; G="Pied droit ";   H=10 },
    @{ Row=502; Name="Omar Benyounes";   C=70; D=3; E=4; F=0; G=$null;           H=0  },
    @{ Row=503; Name="Naim Ighbane";     C=70; D=5; E=0; F=0; G=$null;           H=4  },
    @{ Row=504; Name="Malik Boussaid";   C=70; D=2; E=0; F=0; G=$null;           H=10 },
    @{ Row=505; Name="Romain Thunet";    C=70; D=5; E=5; F=4; G="Genou pizza";   H=7  },
    @{ Row=506; Name="Emmanuel Valey";   C=70; D=6; E=6; F=0; G=$null;           H=7  },
    @{ Row=507; Name="Karahali Souaré";  C=70; D=2; E=6; F=7; G="Ménisque ";     H=0  },
    @{ Row=508; Name="Naim Dhib";        C=70; D=4; E=6; F=3; G="Hanche";        H=4  }
)

foreach ($entry in $data) {
    $r = $entry.Row

    # Reuse formatting from the previous data row (date style, text style,
    # and the special "empty localisation" style) so the new rows look the
    # same as the rest of the table.
    $ws.Range("A498:F498").Copy() | Out-Null
    $ws.Range("A$r`:F$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("H498").Copy() | Out-Null
    $ws.Range("H$r").PasteSpecial(-4122) | Out-Null

    if ($entry.G) {
        $ws.Range("G498").Copy() | Out-Null
        $ws.Range("G$r").PasteSpecial(-4122) | Out-Null
        $ws.Range("G$r").Value = $entry.G
    } else {
        $ws.Range("G496").Copy() | Out-Null
        $ws.Range("G$r").PasteSpecial(-4122) | Out-Null
    }

    $ws.Range("A$r").Value = 45952
    $ws.Range("B$r").Value = $entry.Name
    $ws.Range("C$r").Value = $entry.C
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("E$r").Value = $entry.E
    $ws.Range("F$r").Value = $entry.F
    $ws.Range("H$r").Value = $entry.H
    $ws.Range("I$r").Formula = "=C$r*D$r"
}

$excel.Application.CutCopyMode = $false

$ws.Range("K504").Select()
